$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each changed crypto row

$ws.Range("D2").Value = "96.482.87"
$ws.Range("E2").Value = "  -0.43%  "

$ws.Range("D3").Value = "3.691.63"
$ws.Range("E3").Value = "  +0.56%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.22"
$ws.Range("E5").Value = "  -2.73%  "

$ws.Range("E6").Value = "  -0.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "650.47"
$ws.Range("E7").Value = "  -1.07%  "

$ws.Range("E8").Value = "  +1.19%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -2.57%  "

$ws.Range("D11").Value = "3.688.74"
$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("E12").Value = "  +19.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.26"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("E15").Value = "  +3.27%  "

$ws.Range("D16").Value = "4.378.87"
$ws.Range("E16").Value = "  +0.64%  "

$ws.Range("D17").Value = "96.273.80"
$ws.Range("E17").Value = "  -0.32%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.81"
$ws.Range("E18").Value = "  +9.49%  "

$ws.Range("D19").Value = "3.695.26"
$ws.Range("E19").Value = "  +0.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.92"
$ws.Range("E20").Value = "  -0.07%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.67"
$ws.Range("E21").Value = "  +1.58%  "

$ws.Range("E22").Value = "  -5.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "517.62"
$ws.Range("E23").Value = "  +0.88%  "

$ws.Range("E24").Value = "  -2.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000208"
$ws.Range("E25").Value = "  +1.58%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.94"
$ws.Range("E26").Value = "  +0.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "100.62"
$ws.Range("E27").Value = "  -0.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.11"
$ws.Range("E28").Value = "  +0.23%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.176"
$ws.Range("E29").Value = "  +3.55%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.00"
$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.11"
$ws.Range("E31").Value = "  +1.94%  "

$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.86"
$ws.Range("E33").Value = "  +6.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "653.56"
$ws.Range("E36").Value = "  +6.04%  "

$ws.Range("E37").Value = "  -3.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.587"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.77"
$ws.Range("E39").Value = "  +0.42%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  +12.63%  "

$ws.Range("E42").Value = "  +7.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.38"
$ws.Range("E43").Value = "  -5.69%  "

$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.953"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0449"
$ws.Range("E46").Value = "  +1.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.426"
$ws.Range("E47").Value = "  +1.36%  "

$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("E49").Value = "  -1.78%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.44"
$ws.Range("E50").Value = "  -1.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.52"
$ws.Range("E51").Value = "  +2.34%  "

